# Update the library prep kit name from "E7420" to "E7420L" for all
# sample rows (rows 2-41, column K) to fix inconsistencies found in the
# ZEV jan reps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 11)  # Column K
    if ($cell.Value2 -eq "E7420") {
        $cell.Value2 = "E7420L"
    }
}

# Update the active cell selection to match the author's final cursor
# position in the sheet.
$ws.Range("F35").Select()
